$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update step numbering / text in the "Cenário Normal" flow (rows 13-19)
$ws.Range("C13").Value = "7. Escolhe especificações adicionais"
$ws.Range("D15").Value = "9. Calcula preço"
$ws.Range("D16").Value = "10. Mostra preço e lista de componentes e pergunta se quer confirmar"
$ws.Range("C17").Value = "11. Confirma"
$ws.Range("D18").Value = "12. Regista confirmação"
$ws.Range("D19").Value = "13. Insere carro no sistema"

# Replace "Alternativa 1" block (rows 20-21): now references passo 8 instead of passo 7,
# and the system response text was reworded.
$ws.Range("B20").Value = "Alternativa 1 [Componentes inválidos] (passo 8)"
$ws.Range("D20").Value = "8.1 Verifica que escolheu componentes incompativeis ou estão em falta"
$ws.Range("D21").Value = "8.2 Indica que há peças incompativeis"

# Rows below are unaffected in text but keep them consistent
$ws.Range("D22").Value = "Regressa a 7"
$ws.Range("B24").Value = "Alternativa 2 [Não confirma] (passo 11)"
$ws.Range("C24").Value = "11.1 Não confirma"
$ws.Range("D25").Value = "Regressa a 7"

# Update the view state: scroll position and current selection
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("D21").Select()
